# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (between "2021-Q4" and "总计") holding
#    the per-fund breakdown for the new quarter.
# 2) Prepend a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data (pushing the existing quarters down by one row), and renumber the
#    running index in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the "2022-Q1" sheet by cloning "2021-Q4" so it starts with
# identical layout/styles (header row/text, borders, column-A style, etc.),
# then overwrite the cell values with the new quarter's fund data.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template only has 2 data rows (rows 2-3); the new quarter needs 4
# (rows 2-5), so insert 2 more rows first. Clone the column-A styling
# (running index, s="2") down into the new rows before filling in values.
$newSheet.Rows(4).Insert()
$newSheet.Rows(4).Insert()
$newSheet.Range("A2").Copy($newSheet.Range("A4"))
$newSheet.Range("A2").Copy($newSheet.Range("A5"))

# row -> (基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名)
$fundRows = @(
    @("010714", "东方红远见价值混合", "19.70", "86.34", "4.41", "0.8688", 2),
    @("011351", "金鹰年年邮益一年持有期混合A", "9.03", "37.02", "1.16", "0.1047", 6),
    @("007251", "广发睿享稳健增利混合", "3.69", "38.80", "2.12", "0.0782", 4),
    @("011352", "金鹰年年邮益一年持有期混合C", "0.59", "37.02", "1.16", "0.0068", 6)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    for ($col = 2; $col -le 7; $col++) {
        $cell = $newSheet.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 2]
        $cell.Style = "Normal"
    }

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 2: update the "总计" sheet - insert the new "2022-Q1" row at the top
# of the data (row 2), shifting the existing rows down, then fix up the
# running index in column A for every row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).NumberFormat = "@"
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 2).Style = "Normal"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 3).Style = "Normal"
$summary.Cells.Item(2, 4).Value = 1.06
$summary.Cells.Item(2, 4).Style = "Normal"

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
